$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1810.9062
$ws.Range("I40").Value = 1588.8334
$ws.Range("J40").Value = 2096.4285
$ws.Range("K40").Value = 1588.8334
$ws.Range("L40").Value = 2096.4285
$ws.Range("M40").Value = -1413.8334
$ws.Range("N40").Value = -2446.4285
$ws.Range("H62").Value = 5475.1
$ws.Range("I62").Value = 4761
$ws.Range("J62").Value = 6189.2
$ws.Range("K62").Value = 4761
$ws.Range("L62").Value = 6189.2
$ws.Range("M62").Value = -4137
$ws.Range("N62").Value = -7437.2
$ws.Range("H65").Value = 5475.1
$ws.Range("I65").Value = 4761
$ws.Range("J65").Value = 6189.2
$ws.Range("K65").Value = 23805
$ws.Range("L65").Value = 30946
$ws.Range("M65").Value = -20685
$ws.Range("N65").Value = -37186
$ws.Range("H76").Value = 3165.8538
$ws.Range("I76").Value = 3190.9092
$ws.Range("J76").Value = 3062.5
$ws.Range("K76").Value = 3190.9092
$ws.Range("L76").Value = 3062.5
$ws.Range("M76").Value = -2875.9092
$ws.Range("N76").Value = -3692.5
$ws.Range("H79").Value = 3165.8538
$ws.Range("I79").Value = 3190.9092
$ws.Range("J79").Value = 3062.5
$ws.Range("K79").Value = 3190.9092
$ws.Range("L79").Value = 3062.5
$ws.Range("M79").Value = -2098.9092
$ws.Range("N79").Value = -5246.5
$ws.Range("H80").Value = 4437.34
$ws.Range("I80").Value = 3442.9092
$ws.Range("J80").Value = 5143.0645
$ws.Range("K80").Value = 10328.7276
$ws.Range("L80").Value = 15429.1935
$ws.Range("M80").Value = -9330.7276
$ws.Range("N80").Value = -17425.1935
$ws.Range("H83").Value = 4437.34
$ws.Range("I83").Value = 3442.9092
$ws.Range("J83").Value = 5143.0645
$ws.Range("K83").Value = 30986.1828
$ws.Range("L83").Value = 46287.5805
$ws.Range("M83").Value = -25994.1828
$ws.Range("N83").Value = -56271.5805
$ws.Range("H113").Value = 4970.909
$ws.Range("I113").Value = 2933.3333
$ws.Range("J113").Value = 5735
$ws.Range("K113").Value = 2933.3333
$ws.Range("L113").Value = 5735
$ws.Range("M113").Value = 320.6667000000002
$ws.Range("N113").Value = -12243
$ws.Range("H137").Value = 3139.7778
$ws.Range("I137").Value = 2578.2666
$ws.Range("J137").Value = 3841.6667
$ws.Range("K137").Value = 7734.7998
$ws.Range("L137").Value = 11525.0001
$ws.Range("M137").Value = -5184.7998
$ws.Range("N137").Value = -16625.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1700.1818
$ws.Range("I45").Value = 1204
$ws.Range("J45").Value = 1886.25
$ws.Range("K45").Value = 1204
$ws.Range("L45").Value = 1886.25
$ws.Range("M45").Value = -827
$ws.Range("N45").Value = -2640.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1169.7333
$ws.Range("J80").Value = 404.6
$ws.Range("L80").Value = 404.6
$ws.Range("N80").Value = -2400.6
$ws.Range("H83").Value = 1169.7333
$ws.Range("J83").Value = 404.6
$ws.Range("L83").Value = 2023
$ws.Range("N83").Value = -12007
$ws.Range("H105").Value = 2198.125
$ws.Range("I105").Value = 1597.1428
$ws.Range("J105").Value = 2445.5881
$ws.Range("K105").Value = 1597.1428
$ws.Range("L105").Value = 2445.5881
$ws.Range("M105").Value = 149.8571999999999
$ws.Range("N105").Value = -5939.5881

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 950.4286
$ws.Range("I16").Value = 846.6667
$ws.Range("J16").Value = 1028.25
$ws.Range("K16").Value = 846.6667
$ws.Range("L16").Value = 1028.25
$ws.Range("M16").Value = -559.6667
$ws.Range("N16").Value = -1602.25
$ws.Range("H113").Value = 950.4286
$ws.Range("I113").Value = 846.6667
$ws.Range("J113").Value = 1028.25
$ws.Range("K113").Value = 846.6667
$ws.Range("L113").Value = 1028.25
$ws.Range("M113").Value = 1323.3333
$ws.Range("N113").Value = -5368.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 735.6
$ws.Range("I5").Value = 639.4545000000001
$ws.Range("K5").Value = 1918.3635
$ws.Range("M5").Value = -1806.3635
$ws.Range("H34").Value = 417.44446
$ws.Range("J34").Value = 425.5
$ws.Range("L34").Value = 1276.5
$ws.Range("N34").Value = -1444.5
$ws.Range("H69").Value = 4000
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = -2189
$ws.Range("N69").Value = -16622
$ws.Range("H72").Value = 4000
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 45000
$ws.Range("M72").Value = -4944
$ws.Range("N72").Value = -53112
$ws.Range("H130").Value = 3639.9
$ws.Range("I130").Value = 1399
$ws.Range("J130").Value = 3888.889
$ws.Range("K130").Value = 4197
$ws.Range("L130").Value = 11666.667
$ws.Range("M130").Value = 823
$ws.Range("N130").Value = -21706.667
$ws.Range("H131").Value = 1018.27026
$ws.Range("I131").Value = 826.2222
$ws.Range("J131").Value = 1080
$ws.Range("K131").Value = 2478.6666
$ws.Range("L131").Value = 3240
$ws.Range("M131").Value = 2561.3334
$ws.Range("N131").Value = -13320
$ws.Range("H135").Value = 735.6
$ws.Range("I135").Value = 639.4545000000001
$ws.Range("K135").Value = 5755.0905
$ws.Range("M135").Value = -3220.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 934.2353000000001
$ws.Range("I107").Value = 916.2857
$ws.Range("J107").Value = 946.8
$ws.Range("K107").Value = 916.2857
$ws.Range("L107").Value = 946.8
$ws.Range("M107").Value = 1003.7143
$ws.Range("N107").Value = -4786.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 20222.486
$ws.Range("I93").Value = 1449.862
$ws.Range("J93").Value = 88273.25
$ws.Range("K93").Value = 1449.862
$ws.Range("L93").Value = 88273.25
$ws.Range("M93").Value = -201.8620000000001
$ws.Range("N93").Value = -90769.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 518.125
$ws.Range("I113").Value = 375
$ws.Range("J113").Value = 947.5
$ws.Range("K113").Value = 1125
$ws.Range("L113").Value = 2842.5
$ws.Range("M113").Value = 1045
$ws.Range("N113").Value = -7182.5
